# Auto-generated: apply cryptos price/volume update (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force column D price strings (e.g. "1.00", "3.168.14") to stay
    # plain text instead of being auto-coerced to a number/date by Excel,
    # then restore the cell style so no stray formatting is introduced.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "79.165.60"
$ws.Range("E2").Value = "  +3.60%  "

Set-TextValue $ws.Range("D3") "3.184.58"
$ws.Range("E3").Value = "  +7.07%  "

Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  -0.04%  "

Set-TextValue $ws.Range("D5") "206.78"
$ws.Range("E5").Value = "  +3.22%  "

Set-TextValue $ws.Range("D6") "633.61"
$ws.Range("E6").Value = "  +0.38%  "

Set-TextValue $ws.Range("D7") "0.999"
$ws.Range("E7").Value = "  -0.02%  "

Set-TextValue $ws.Range("D8") "0.229"
$ws.Range("E8").Value = "  +13.51%  "

Set-TextValue $ws.Range("D9") "0.577"
$ws.Range("E9").Value = "  +5.26%  "

Set-TextValue $ws.Range("D10") "3.183.82"
$ws.Range("E10").Value = "  +7.10%  "

Set-TextValue $ws.Range("D11") "0.567"
$ws.Range("E11").Value = "  +32.41%  "

Set-TextValue $ws.Range("D12") "0.163"
$ws.Range("E12").Value = "  +1.41%  "

Set-TextValue $ws.Range("D13") "5.44"
$ws.Range("E13").Value = "  +9.46%  "

Set-TextValue $ws.Range("D14") "3.762.22"
$ws.Range("E14").Value = "  +7.04%  "

Set-TextValue $ws.Range("D15") "0.0000228"
$ws.Range("E15").Value = "  +21.48%  "

Set-TextValue $ws.Range("D16") "31.61"
$ws.Range("E16").Value = "  +9.15%  "

Set-TextValue $ws.Range("D17") "78.828.01"
$ws.Range("E17").Value = "  +3.36%  "

Set-TextValue $ws.Range("D18") "3.181.35"
$ws.Range("E18").Value = "  +7.40%  "

Set-TextValue $ws.Range("D19") "14.35"
$ws.Range("E19").Value = "  +7.17%  "

Set-TextValue $ws.Range("D20") "9.43"
$ws.Range("E20").Value = "  +6.57%  "

Set-TextValue $ws.Range("D21") "432.80"
$ws.Range("E21").Value = "  +16.44%  "

Set-TextValue $ws.Range("D22") "2.84"
$ws.Range("E22").Value = "  +25.45%  "

Set-TextValue $ws.Range("D23") "4.88"
$ws.Range("E23").Value = "  +13.83%  "

Set-TextValue $ws.Range("D24") "6.83"
$ws.Range("E24").Value = "  +6.50%  "

Set-TextValue $ws.Range("D29") "1.01"
$ws.Range("E29").Value = "  +0.96%  "

Set-TextValue $ws.Range("D30") "0.0000116"
$ws.Range("E30").Value = "  +9.09%  "

Set-TextValue $ws.Range("D31") "0.996"
$ws.Range("E31").Value = "  -0.44%  "

Set-TextValue $ws.Range("D32") "8.88"
$ws.Range("E32").Value = "  +7.82%  "

Set-TextValue $ws.Range("D33") "1.47"
$ws.Range("E33").Value = "  +6.44%  "

Set-TextValue $ws.Range("D34") "519.14"
$ws.Range("E34").Value = "  +2.27%  "

$ws.Range("E35").Value = "  +1.70%  "

$ws.Range("E36").Value = "  +22.67%  "

Set-TextValue $ws.Range("D37") "22.51"
$ws.Range("E37").Value = "  +11.02%  "

Set-TextValue $ws.Range("D38") "1.00"
$ws.Range("E38").Value = "  +0.01%  "

Set-TextValue $ws.Range("D39") "0.398"
$ws.Range("E39").Value = "  +4.51%  "

Set-TextValue $ws.Range("D40") "164.20"
$ws.Range("E40").Value = "  +0.15%  "

Set-TextValue $ws.Range("D41") "197.09"
$ws.Range("E41").Value = "  +6.28%  "

$ws.Range("E44").Value = "  -0.24%  "

Set-TextValue $ws.Range("D45") "5.42"
$ws.Range("E45").Value = "  +10.15%  "

Set-TextValue $ws.Range("D46") "0.801"
$ws.Range("E46").Value = "  +14.92%  "

Set-TextValue $ws.Range("D47") "1.79"
$ws.Range("E47").Value = "  +9.41%  "

$ws.Range("E48").Value = "  +5.55%  "

Set-TextValue $ws.Range("D51") "0.625"
$ws.Range("E51").Value = "  +6.59%  "

# Full row updates (coins re-ranked by the scraper)
$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D25") "4.79"
$ws.Range("E25").Value = "  +10.80%  "

$ws.Range("B26").Value = "Aptos"
$ws.Range("C26").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D26") "11.15"
$ws.Range("E26").Value = "  +14.73%  "

$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue $ws.Range("D27") "3.335.54"
$ws.Range("E27").Value = "  +6.61%  "

$ws.Range("B28").Value = "Litecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D28") "76.57"
$ws.Range("E28").Value = "  +5.19%  "

$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D42") "0.109"
$ws.Range("E42").Value = "  +4.12%  "

$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue $ws.Range("D43") "20.00"
$ws.Range("E43").Value = "  +0.10%  "

$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D49") "2.58"
$ws.Range("E49").Value = "  +11.93%  "

$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D50") "42.88"
$ws.Range("E50").Value = "  +0.49%  "

